$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 151.125
$ws.Range("I33").Value = 77.05556
$ws.Range("J33").Value = 373.33334
$ws.Range("K33").Value = 77.05556
$ws.Range("L33").Value = 373.33334
$ws.Range("M33").Value = 151.94444
$ws.Range("N33").Value = -831.33334
# Row 64
$ws.Range("H64").Value = 3812.5942
$ws.Range("I64").Value = 3711.628
$ws.Range("J64").Value = 3979.577
$ws.Range("K64").Value = 3711.628
$ws.Range("L64").Value = 3979.577
$ws.Range("M64").Value = -3463.628
$ws.Range("N64").Value = -4475.577
# Row 67
$ws.Range("H67").Value = 3812.5942
$ws.Range("I67").Value = 3711.628
$ws.Range("J67").Value = 3979.577
$ws.Range("K67").Value = 3711.628
$ws.Range("L67").Value = 3979.577
$ws.Range("M67").Value = -2853.628
$ws.Range("N67").Value = -5695.577
# Row 101
$ws.Range("H101").Value = 408.625
$ws.Range("I101").Value = 352.7143
$ws.Range("J101").Value = 800
$ws.Range("K101").Value = 1058.1429
$ws.Range("L101").Value = 2400
$ws.Range("M101").Value = 563.8571000000002
$ws.Range("N101").Value = -5644

$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Range("H122").Value = 791.8570999999999
$ws.Range("I122").Value = 791.8570999999999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2375.5713
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 74.42870000000039
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1462.3103
$ws.Range("I134").Value = 1346.25
$ws.Range("K134").Value = 4038.75
$ws.Range("M134").Value = -1503.75

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 17862.5
$ws.Range("I17").Value = 450
$ws.Range("J17").Value = 23666.666
$ws.Range("K17").Value = 450
$ws.Range("L17").Value = 23666.666
$ws.Range("M17").Value = -276
$ws.Range("N17").Value = -24014.666
# Row 51
$ws.Range("H51").Value = 25000
$ws.Range("J51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -26472
# Row 58
$ws.Range("H58").Value = 2016.7812
$ws.Range("I58").Value = 1486.619
$ws.Range("J58").Value = 3028.9092
$ws.Range("K58").Value = 1486.619
$ws.Range("L58").Value = 3028.9092
$ws.Range("M58").Value = -1283.619
$ws.Range("N58").Value = -3434.9092
# Row 60
$ws.Range("H60").Value = 10802.608
$ws.Range("I60").Value = 9500
$ws.Range("J60").Value = 10861.818
$ws.Range("K60").Value = 9500
$ws.Range("L60").Value = 10861.818
$ws.Range("M60").Value = -8989
$ws.Range("N60").Value = -11883.818
# Row 61
$ws.Range("H61").Value = 25000
$ws.Range("J61").Value = 25000
$ws.Range("L61").Value = 25000
$ws.Range("N61").Value = -25696
# Row 62
$ws.Range("H62").Value = 2001659.6
$ws.Range("I62").Value = 2001659.6
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2001659.6
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2001035.6
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 2001659.6
$ws.Range("I65").Value = 2001659.6
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 10008298
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -10005178
$ws.Range("N65").ClearContents()
# Row 136
$ws.Range("H136").Value = 2016.7812
$ws.Range("I136").Value = 1486.619
$ws.Range("J136").Value = 3028.9092
$ws.Range("K136").Value = 4459.857
$ws.Range("L136").Value = 9086.7276
$ws.Range("M136").Value = -1909.857
$ws.Range("N136").Value = -14186.7276

$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Range("H9").Value = 54350
$ws.Range("I9").Value = 167333.33
$ws.Range("J9").Value = 5928.5713
$ws.Range("K9").Value = 501999.99
$ws.Range("L9").Value = 17785.7139
$ws.Range("M9").Value = -501775.99
$ws.Range("N9").Value = -18233.7139
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
# Row 35
$ws.Range("H35").Value = 1053.0667
$ws.Range("I35").Value = 1200
$ws.Range("J35").Value = 1030.4615
$ws.Range("K35").Value = 3600
$ws.Range("L35").Value = 3091.3845
$ws.Range("M35").Value = -3312
$ws.Range("N35").Value = -3667.3845
# Row 49
$ws.Range("H49").Value = 3001.3333
$ws.Range("J49").Value = 3001.3333
$ws.Range("L49").Value = 9003.999899999999
$ws.Range("N49").Value = -9315.999899999999
# Row 54
$ws.Range("H54").Value = 2000
$ws.Range("I54").Value = 2000
$ws.Range("K54").Value = 6000
$ws.Range("M54").Value = -5441
# Row 57
$ws.Range("H57").Value = 4200
$ws.Range("J57").Value = 5000
$ws.Range("L57").Value = 15000
$ws.Range("N57").Value = -16118
# Row 101
$ws.Range("H101").Value = 4803.222
$ws.Range("J101").Value = 4928.625
$ws.Range("L101").Value = 14785.875
$ws.Range("N101").Value = -19653.875
# Row 112
$ws.Range("H112").Value = 3068.4
$ws.Range("I112").Value = 1531.5
$ws.Range("J112").Value = 3627.2727
$ws.Range("K112").Value = 4594.5
$ws.Range("L112").Value = 10881.8181
$ws.Range("M112").Value = -3486.5
$ws.Range("N112").Value = -13097.8181
# Row 113
$ws.Range("H113").Value = 2029105.9
$ws.Range("I113").Value = 5747700
$ws.Range("J113").Value = 781.8182
$ws.Range("K113").Value = 17243100
$ws.Range("L113").Value = 2345.4546
$ws.Range("M113").Value = -17240930
$ws.Range("N113").Value = -6685.4546
# Row 118
$ws.Range("H118").Value = 2411.6365
$ws.Range("I118").Value = 609.6667
$ws.Range("J118").Value = 3087.375
$ws.Range("K118").Value = 1829.0001
$ws.Range("L118").Value = 9262.125
$ws.Range("M118").Value = -586.0001
$ws.Range("N118").Value = -11748.125

$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 20000
$ws.Range("J62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("N62").Value = -21372
# Row 65
$ws.Range("H65").Value = 20000
$ws.Range("J65").Value = 20000
$ws.Range("L65").Value = 60000
$ws.Range("N65").Value = -66864
# Row 136
$ws.Range("H136").Value = 24660.5
$ws.Range("J136").Value = 24660.5
$ws.Range("L136").Value = 73981.5
$ws.Range("N136").Value = -79081.5

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 6921.1
$ws.Range("I40").Value = 7173
$ws.Range("J40").Value = 6333.3335
$ws.Range("K40").Value = 7173
$ws.Range("L40").Value = 6333.3335
$ws.Range("M40").Value = -7037
$ws.Range("N40").Value = -6605.3335
# Row 46
$ws.Range("H46").Value = 1133.3334
$ws.Range("I46").Value = 700
$ws.Range("K46").Value = 700
$ws.Range("M46").Value = -512
# Row 136
$ws.Range("H136").Value = 1687.6111
$ws.Range("I136").Value = 994.38464
$ws.Range("J136").Value = 3490
$ws.Range("K136").Value = 2983.15392
$ws.Range("L136").Value = 10470
$ws.Range("M136").Value = -433.1539199999997
$ws.Range("N136").Value = -15570

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1350
$ws.Range("I81").Value = 916.6667
$ws.Range("K81").Value = 1833.3334
$ws.Range("M81").Value = -772.3334
# Row 84
$ws.Range("H84").Value = 1350
$ws.Range("I84").Value = 916.6667
$ws.Range("K84").Value = 9166.666999999999
$ws.Range("M84").Value = -3862.666999999999
# Row 113
$ws.Range("H113").Value = 338.9091
$ws.Range("I113").Value = 333
$ws.Range("K113").Value = 999
$ws.Range("M113").Value = 1171
